$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "updated activity till excel form" - the per-innings batting stats
# (runs/balls/fours/sixes) for rows 2,4,5,6,7,8 are refreshed with the
# latest figures. Row 3 is untouched. All of these columns are stored
# as text (see the sheet's numberStoredAsText ignored-error range), so
# each write is wrapped with a Text number format to keep the cell a
# string instead of letting Excel auto-convert the digits to a number,
# then the format is reset back to Normal/General so no stray style is
# left behind on the cell.

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "C2" "70"
Set-TextValue "D2" "44"
Set-TextValue "E2" "4"
Set-TextValue "F2" "5"

Set-TextValue "C4" "16"
Set-TextValue "D4" "13"
Set-TextValue "E4" "1"
Set-TextValue "F4" "1"

Set-TextValue "C5" "22"
Set-TextValue "D5" "12"
Set-TextValue "E5" "3"
Set-TextValue "F5" "1"

Set-TextValue "C6" "4"
Set-TextValue "D6" "7"
Set-TextValue "E6" "0"
Set-TextValue "F6" "0"

Set-TextValue "C7" "21"
Set-TextValue "D7" "16"
Set-TextValue "E7" "1"
Set-TextValue "F7" "2"

Set-TextValue "C8" "13"
Set-TextValue "D8" "8"
Set-TextValue "E8" "2"
Set-TextValue "F8" "0"
